$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R ("2021") added next to the existing year columns (D..Q = 2007..2020).
# Each target cell's number format/font/border is reproduced by copying the
# formatting from an existing cell in the sheet that already carries the
# desired style (header year cell, bold region/total row, or plain data row),
# then overwriting just the value - this reuses the workbook's existing
# style entries instead of minting unnecessary new ones.

$rows = @(
  @{Row=4;  Src='Q4';  Val=2021},
  @{Row=5;  Src='A5';  Val=5.3},
  @{Row=6;  Src='A10'; Val=6.3},
  @{Row=7;  Src='A10'; Val=4.7},
  @{Row=8;  Src='A8';  Val=$null},
  @{Row=9;  Src='A5';  Val=6.6},
  @{Row=10; Src='A10'; Val=7.5},
  @{Row=11; Src='A10'; Val=6.2},
  @{Row=12; Src='A5';  Val=11.8},
  @{Row=13; Src='A10'; Val=15.5},
  @{Row=14; Src='A10'; Val=9.6999999999999993},
  @{Row=15; Src='A5';  Val=6.3},
  @{Row=16; Src='A10'; Val=7.5},
  @{Row=17; Src='A10'; Val=5.6},
  @{Row=18; Src='A5';  Val=6.3},
  @{Row=19; Src='A10'; Val=10.8},
  @{Row=20; Src='A10'; Val=4.3},
  @{Row=21; Src='A5';  Val=1.9},
  @{Row=22; Src='A10'; Val=3.1},
  @{Row=23; Src='A10'; Val=1.1000000000000001},
  @{Row=24; Src='A5';  Val=2.6},
  @{Row=25; Src='A10'; Val=3.8},
  @{Row=26; Src='A10'; Val=1.7},
  @{Row=27; Src='A5';  Val=5.3},
  @{Row=28; Src='A10'; Val=6.2},
  @{Row=29; Src='A10'; Val=4.8},
  @{Row=30; Src='A5';  Val=4.0999999999999996},
  @{Row=31; Src='A10'; Val=3.3},
  @{Row=32; Src='A10'; Val=4.9000000000000004},
  @{Row=33; Src='A5';  Val=2.8},
  @{Row=34; Src='A10'; Val=3.4},
  @{Row=35; Src='A10'; Val=2.6},
  @{Row=36; Src='A8';  Val=$null},
  @{Row=37; Src='A10'; Val=15.7},
  @{Row=38; Src='A10'; Val=7.9},
  @{Row=39; Src='A10'; Val=4.5},
  @{Row=40; Src='A10'; Val=4.4000000000000004},
  @{Row=41; Src='A10'; Val=2.9},
  @{Row=42; Src='A10'; Val=1.4}
)

foreach ($item in $rows) {
  $target = $ws.Range("R$($item.Row)")
  $ws.Range($item.Src).Copy()
  $target.PasteSpecial(-4122)
  if ($null -ne $item.Val) {
    $target.Value = $item.Val
  }
}

# Row 43 ("70 and over") gets the bottom-border text style (copied from A43)
# with right alignment, holding the same ellipsis "..." shared string already
# used elsewhere in the row (D43, E43, F43, I43, O43, P43 = shared-string 7).
$ws.Range("A43").Copy()
$r43 = $ws.Range("R43")
$r43.PasteSpecial(-4122)
$r43.HorizontalAlignment = -4152
$r43.Value = "…"

$excel.CutCopyMode = $false

$ws.Range("S1").Select()
"done"
